$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 184-185, pushing the existing rows 184-189
# down to 186-191 (new weekly price data is prepended, most-recent-first).
$ws.Range("A184:A185").EntireRow.Insert()

# --- New row 184 (Primera) ---
$ws.Cells.Item(184, 1).Value = 11
$ws.Cells.Item(184, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(184, 3).Value = "Bíobío"
$ws.Cells.Item(184, 4).Value = 44610
$ws.Cells.Item(184, 5).Value = 8
$ws.Cells.Item(184, 6).Value = 100114013
$ws.Cells.Item(184, 7).Value = "Zanahoria"
$ws.Cells.Item(184, 8).Value = "Sin especificar"
$ws.Cells.Item(184, 9).Value = "Primera"
$ws.Cells.Item(184, 10).Value = 800
$ws.Cells.Item(184, 11).Value = 8000
$ws.Cells.Item(184, 12).Value = 8500
$ws.Cells.Item(184, 13).Value = 8250
$ws.Cells.Item(184, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(184, 15).Value = "Región de Ñuble"
$ws.Cells.Item(184, 16).Value = 412
$ws.Cells.Item(184, 17).Value = 20
$ws.Cells.Item(184, 18).Value = "Hortaliza"

# --- New row 185 (Segunda) ---
$ws.Cells.Item(185, 1).Value = 11
$ws.Cells.Item(185, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(185, 3).Value = "Bíobío"
$ws.Cells.Item(185, 4).Value = 44610
$ws.Cells.Item(185, 5).Value = 8
$ws.Cells.Item(185, 6).Value = 100114013
$ws.Cells.Item(185, 7).Value = "Zanahoria"
$ws.Cells.Item(185, 8).Value = "Sin especificar"
$ws.Cells.Item(185, 9).Value = "Segunda"
$ws.Cells.Item(185, 10).Value = 400
$ws.Cells.Item(185, 11).Value = 7000
$ws.Cells.Item(185, 12).Value = 7000
$ws.Cells.Item(185, 13).Value = 7000
$ws.Cells.Item(185, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(185, 15).Value = "Región de Ñuble"
$ws.Cells.Item(185, 16).Value = 350
$ws.Cells.Item(185, 17).Value = 20
$ws.Cells.Item(185, 18).Value = "Hortaliza"
